# Add the new "Calls" worksheet (Word Mail Merge Template data) after the
# existing "Cases" sheet, matching the freeCrmTestData workbook layout.

$wb = $excel.ActiveWorkbook

$casesSheet = $wb.Worksheets.Item("Cases")
$newSheet = $wb.Worksheets.Add($null, $casesSheet)
$newSheet.Name = "Calls"

# Header row (row 1) - give it the same yellow-fill header style used by
# the other sheets (reuses the workbook's existing header style).
$newSheet.Range("A1:E1").Interior.Color = 65535
$newSheet.Range("A1").Value = "contact"
$newSheet.Range("B1").Value = "deal"
$newSheet.Range("C1").Value = "task"
$newSheet.Range("D1").Value = "case"
$newSheet.Range("E1").Value = "notes"

# Column A data
$newSheet.Range("A2").Value = "zzzx"
$newSheet.Range("A3").Value = "aaaa"

# Column B data
$newSheet.Range("B2").Value = "aaaa"
$newSheet.Range("B3").Value = "bbbb"

# Column C data
$newSheet.Range("C2").Value = "cccc"
$newSheet.Range("C3").Value = "dddd"

# Column D data
$newSheet.Range("D2").Value = "eeee"
$newSheet.Range("D3").Value = "ffff"

# Column E data
$newSheet.Range("E2").Value = "gggg"
$newSheet.Range("E3").Value = "hhhh"

# Match the recorded UI state: column B selected on the new Calls sheet ...
$newSheet.Columns.Item(2).Select() | Out-Null

# ... and the whole first row selected on the Cases sheet.
$casesSheet.Rows.Item(1).Select() | Out-Null

# Calls ends up the active/selected tab.
$newSheet.Activate()
